# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the default "Table_0" style ({749AE586-0691-418A-882A-76309F03B465})
#    to the built-in style {A50A1380-744F-40B0-8E7B-CE56D5281469}.
#
# 2) The deck's theme is swapped: the master/theme that drives the slides'
#    look ("Integral" / "Red Violet" colour scheme) is switched for the
#    plain default "Office" colour scheme (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), while the font scheme and format scheme (already
#    identical between the two themes in this deck) are left untouched.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------

$newTableStyleId = "{A50A1380-744F-40B0-8E7B-CE56D5281469}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# --- 2) Swap the presentation's colour theme -------------------------------
# The deck's slide master currently uses the "Integral" theme (Red Violet
# colour scheme); flip it to the standard Office colour scheme instead
# (the font scheme / format scheme already match, only the 12 theme
# colours differ).

$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
